$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear stale cells whose content moved to a new location ---
$ws.Range("A27").ClearContents()
$ws.Range("C27").ClearContents()

# --- Pre-existing cell values (reuse existing shared strings 0-29) ---
$ws.Range("B1").Value = 'Hillz'
$ws.Range("C1").Value = 'Trump'
$ws.Range("A2").Value = 'Book'
$ws.Range("B2").Value = 'Fyodor Dostoyevsky''s ''The Brothers Karamazov'''
$ws.Range("C2").Value = 'All Quiet on the Western Front'
$ws.Range("B3").Value = 'West with the Night'
$ws.Range("C3").Value = 'The Power of Positive Thinking'
$ws.Range("C4").Value = 'Essays and Lectures by Ralph Waldo Emerson'
$ws.Range("B8").Value = 'I love to swim'
$ws.Range("A10").Value = 'Color'
$ws.Range("B10").Value = 'Yellow'
$ws.Range("C10").Value = 'Gold'
$ws.Range("A11").Value = 'Snack'
$ws.Range("B11").Value = 'Chocolate'
$ws.Range("C11").Value = 'Bacon, Egg, and Toast Cups'
$ws.Range("A12").Value = 'Snack'
$ws.Range("B12").Value = 'Fruit'
$ws.Range("A15").Value = 'Movie'
$ws.Range("B15").Value = 'Wizard of Oz'
$ws.Range("C15").Value = 'Citizen Kane'
$ws.Range("A16").Value = 'Movie'
$ws.Range("B16").Value = 'Casablanca'
$ws.Range("A18").Value = 'birthplace'
$ws.Range("B18").Value = 'Edgewater Hospital, Chicago'
$ws.Range("C18").Value = 'Queens, New York City, NY'
$ws.Range("A20").Value = 'birthday'
$ws.Range("B20").Value = 'October 26, 1947 (age 68 years)'
$ws.Range("C20").Value = 'June 14, 1946 (age 70 years)'
$ws.Range("B29").Value = 'president of her high school class'
$ws.Range("A31").Value = 'Pets'
$ws.Range("B31").Value = 'Buddy (dog) Buddy (August 7, 1997 – January 2, 2002), a male chocolate-colored Labrador Retriever, was one of two pets kept by the Clinton family while Bill Clinton was President of the United States. The Clintons'' other pet was a cat named Socks.'
$ws.Range("C31").Value = 'a dog named Spinee???'

# --- B22: date value (April 12, 2015) with date number format (style index 2) ---
$ws.Range("B22").NumberFormat = "d-mmm-yy"
$ws.Range("B22").Value = "4/12/2015"

# --- New cell values, in the exact order the source workbook registered them ---
# --- (so the generated shared-string table indices line up with the target) ---
$ws.Range("A22").Value = 'announced their bid for presidency'
$ws.Range("B24").Value = '67th United States Secretary of State'
$ws.Range("B25").Value = 'YouTube video'
$ws.Range("A25").Value = 'announced via '
$ws.Range("A33").Value = 'Hillary was in the 2015 Time magazine''s "100 Most Influential People". Which of these people were also in that list?'
$ws.Range("B33").Value = 'Amy Schumer, Bradley Cooper, Kim Kardashian West, Ina Garten'
$ws.Range("A34").Value = 'Office'
$ws.Range("B34").Value = 'Clinton had taken a lease on a small office at 1 Pierrepont Plaza in Brooklyn, New York City -  Morgan Stanley has a major office in the building, which is also the home of the law office of Loretta E. Lynch'
$ws.Range("B26").Value = 'She stated that, "Everyday Americans need a champion. And I want to be that champion.'
$ws.Range("B23").NumberFormat = "d-mmm-yy"
$ws.Range("B23").WrapText = $true
$ws.Range("B23").Value = 'while Florida Senator Marco Rubio announced his candidacy on April 13, the day after Clinton. '
$ws.Range("B27").Value = 'Immediately following her announcement, she made a two-day road trip in a customized Chevrolet Express van, nicknamed after Scooby-Doo'
$ws.Range("B28").Value = 'Clinton held her first major campaign rally June 13, 2015, at Franklin D. Roosevelt Four Freedoms Park on the southern tip of New York City''s Roosevelt Island.'
$ws.Range("B36").Value = 'Graphic designer Michael Bierut of the firm Pentagram designed the campaign''s distinctive "H" logo;'
$ws.Range("B37").Value = 'Headquarters Brooklyn, New York, U.S.'
$ws.Range("B39").Value = 'Newsweek ranked her as the 13th most powerful person on the planet'
$ws.Range("B40").Value = 'In 2012, she was chosen as one of Barbara Walters'' 10 Most Fascinating People of the year.[40]'
$ws.Range("B42").Value = 'Clinton has been ranked on their list of the world''s most powerful people by Forbes magazine. She was listed as 5th most powerful in 2004,[67] 26th in 2005,[68] 18th in 2006,[69] 28th in 2008,[70] 36th in 2009,[71] 2nd in 2011,[72] 2nd in 2012,[73] 5th in 2013,[74] 6th in 2014, and 58th in 2015.[75][76]'
$ws.Range("B44").Value = 'Clinton has been named ten times in Time magazine''s Time 100 as one of the 100 most influential people in the world.[77][78][79] Years this happened were 2004 (as part of The Clintons),[80] 2006,[81] 2007,[82] 2008,[83] 2009,[84] 2011,[85] 2012,[86] 2014,[77] 2015,[78] and 2016.[79] In addition, in November 2010, Time named Clinton one of the 25 most powerful women of the past century.[87]'
$ws.Range("B45").Value = 'Clinton has been named three times as Barbara Walters'' Most Fascinating Person of the year, in 1993, 2003, and 2013.[88]'
$ws.Range("B47").Value = 'Family 
Bill Clinton (Husband presidency) Chelsea Clinton (Daughter) Hugh E. Rodham (Father) Dorothy Howell Rodham (Mother) Hugh Rodham (Brother) Tony Rodham (Brother) Socks (Cat) Buddy (Dog) Whitehaven (house)'
$ws.Range("B48").Value = 'Writings 
Bibliography Senior thesis (1969) It Takes a Village (1996) Dear Socks, Dear Buddy (1998) An Invitation to the White House (2000) Living History (2003) Hard Choices (2014)'
$ws.Range("A50").Value = 'Democratic blue to Republican red. A purple state refers to a swing state where both Democratic and Republican candidates receive strong support without an overwhelming majority of support for either party.'

# --- Wrap text for all other populated cells (style index 1) ---
$ws.Range("B1,C1,A2,B2,C2,B3,C3,C4,B8,A10,B10,C10,A11,B11,C11,A12,B12,A15,B15,C15,A16,B16,A18,B18,C18,A20,B20,C20,A22,B24,A25,B25,B26,B27,B28,B29,A31,B31,C31,A33,B33,A34,B34,B36,B37,B39,B40,B42,B44,B45,B47,B48,A50").WrapText = $true

# --- Row heights ---
$ws.Rows.Item(23).RowHeight = 32
$ws.Rows.Item(26).RowHeight = 32
$ws.Rows.Item(27).RowHeight = 32
$ws.Rows.Item(28).RowHeight = 48
$ws.Rows.Item(31).RowHeight = 64
$ws.Rows.Item(33).RowHeight = 48
$ws.Rows.Item(34).RowHeight = 48
$ws.Rows.Item(36).RowHeight = 32
$ws.Rows.Item(40).RowHeight = 32
$ws.Rows.Item(42).RowHeight = 64
$ws.Rows.Item(44).RowHeight = 96
$ws.Rows.Item(45).RowHeight = 32
$ws.Rows.Item(47).RowHeight = 64
$ws.Rows.Item(48).RowHeight = 64
$ws.Rows.Item(50).RowHeight = 64

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 44.83203125
$ws.Columns.Item(2).ColumnWidth = 65.66796875

# --- View: selection + zoom ---
$ws.Range("A33").Select()
$excel.ActiveWindow.Zoom = 105
